$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "2023-12-07 13:01:54"
$ws.Range("B69").Value = 0.0002
